# Fruta / hortaliza, semanal
# A new weekly price-report row is published for "Terminal Hortofrutícola
# Agro Chillán - Mango": it lands at the top of the data table (row 27,
# right under the header-less block start) and every existing record from
# row 27 down gets pushed one row further (27->28, 28->29, ..., 44->45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 27, shifting rows 27:44 down to 28:45 (cells/format
# move with their row, same as Excel's native Rows(n).Insert()).
$ws.Rows(27).Insert()

# Populate the newly opened row 27 with this week's record.
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = 44447
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108002
$ws.Range("J27").Value = "Mango"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 60
$ws.Range("N27").Value = 8500
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 8750
$ws.Range("Q27").Value = "$/bandeja 4 kilos"
$ws.Range("R27").Value = "Brasil"
$ws.Range("S27").Value = 2188
$ws.Range("T27").Value = 4
